$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Apply a thin border around every cell in the table (A1:D34), and give
#    the header row (A1:D1) bold + centered text on top of that border.
#    Build each look on a scratch cell first and paste the *formats only*
#    onto the target range in a single operation - this keeps the resulting
#    styles.xml compact (one cellXf per distinct resulting format), matching
#    how the workbook was actually edited.
# ---------------------------------------------------------------------------
$tmplBody = $ws.Range("Z1")
$tmplBody.Borders.LineStyle = 1
$tmplBody.Copy()
$ws.Range("A2:D34").PasteSpecial(-4122)   # xlPasteFormats
$tmplBody.Clear()

$tmplHeader = $ws.Range("Z2")
$tmplHeader.Font.Bold = $true
$tmplHeader.HorizontalAlignment = -4108   # xlCenter
$tmplHeader.Borders.LineStyle = 1
$tmplHeader.Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)    # xlPasteFormats
$tmplHeader.Clear()

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Widen column C.
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 37

# ---------------------------------------------------------------------------
# 3. Update the view: drop the old scroll position / selection and select
#    the whole of row 9 instead.
# ---------------------------------------------------------------------------
$ws.Rows(9).Select()

# ---------------------------------------------------------------------------
# 4. Page setup: paper size A4 (9) and portrait orientation.
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
